$d = $word.ActiveDocument

# Locate the paragraph that contains the "© 2020 ..." footer text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Powered by Jekyll and Github pages*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    # The two empty paragraphs immediately preceding the footer paragraph
    # (one plain, one with a page-break-before) are removed along with the
    # footer paragraph itself.
    $startPara = $d.Paragraphs.Item($target - 2)
    $endPara = $d.Paragraphs.Item($target)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
